$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("existing_stock")

# Update column P (description) for rows 64-88 with reshuffled
# "Aggregated Plant - IRENA Gap - CHE_n_Missing Solar Capacity" labels
$ws.Range("P64").Value = "Aggregated Plant - IRENA Gap - CHE_19_Missing Solar Capacity"
$ws.Range("P65").Value = "Aggregated Plant - IRENA Gap - CHE_14_Missing Solar Capacity"
$ws.Range("P66").Value = "Aggregated Plant - IRENA Gap - CHE_10_Missing Solar Capacity"
$ws.Range("P67").Value = "Aggregated Plant - IRENA Gap - CHE_0_Missing Solar Capacity"
$ws.Range("P68").Value = "Aggregated Plant - IRENA Gap - CHE_1_Missing Solar Capacity"
$ws.Range("P69").Value = "Aggregated Plant - IRENA Gap - CHE_24_Missing Solar Capacity"
$ws.Range("P70").Value = "Aggregated Plant - IRENA Gap - CHE_23_Missing Solar Capacity"
$ws.Range("P71").Value = "Aggregated Plant - IRENA Gap - CHE_6_Missing Solar Capacity"
$ws.Range("P72").Value = "Aggregated Plant - IRENA Gap - CHE_17_Missing Solar Capacity"
$ws.Range("P73").Value = "Aggregated Plant - IRENA Gap - CHE_3_Missing Solar Capacity"
$ws.Range("P74").Value = "Aggregated Plant - IRENA Gap - CHE_20_Missing Solar Capacity"
$ws.Range("P75").Value = "Aggregated Plant - IRENA Gap - CHE_25_Missing Solar Capacity"
$ws.Range("P76").Value = "Aggregated Plant - IRENA Gap - CHE_8_Missing Solar Capacity"
$ws.Range("P77").Value = "Aggregated Plant - IRENA Gap - CHE_5_Missing Solar Capacity"
$ws.Range("P78").Value = "Aggregated Plant - IRENA Gap - CHE_22_Missing Solar Capacity"
$ws.Range("P79").Value = "Aggregated Plant - IRENA Gap - CHE_9_Missing Solar Capacity"
$ws.Range("P80").Value = "Aggregated Plant - IRENA Gap - CHE_18_Missing Solar Capacity"
$ws.Range("P81").Value = "Aggregated Plant - IRENA Gap - CHE_2_Missing Solar Capacity"
$ws.Range("P82").Value = "Aggregated Plant - IRENA Gap - CHE_4_Missing Solar Capacity"
$ws.Range("P83").Value = "Aggregated Plant - IRENA Gap - CHE_11_Missing Solar Capacity"
$ws.Range("P84").Value = "Aggregated Plant - IRENA Gap - CHE_12_Missing Solar Capacity"
$ws.Range("P85").Value = "Aggregated Plant - IRENA Gap - CHE_7_Missing Solar Capacity"
$ws.Range("P86").Value = "Aggregated Plant - IRENA Gap - CHE_21_Missing Solar Capacity"
$ws.Range("P87").Value = "Aggregated Plant - IRENA Gap - CHE_13_Missing Solar Capacity"
$ws.Range("P88").Value = "Aggregated Plant - IRENA Gap - CHE_15_Missing Solar Capacity"

# Update column E (ncap_pasti) for rows 99-123 with reshuffled random values
$ws.Range("E99").Value = 0.1663793251004252
$ws.Range("E100").Value = 0.1727926292604506
$ws.Range("E101").Value = 0.16277902359433066
$ws.Range("E102").Value = 0.15291072157643879
$ws.Range("E103").Value = 0.16554834018408843
$ws.Range("E104").Value = 0.21471510601685545
$ws.Range("E105").Value = 0.15409499379434963
$ws.Range("E106").Value = 0.16152827258311295
$ws.Range("E107").Value = 0.1371001651339535
$ws.Range("E108").Value = 0.18308354646436523
$ws.Range("E109").Value = 0.15521278721895346
$ws.Range("E110").Value = 0.1586174511333161
$ws.Range("E111").Value = 0.20091315882928704
$ws.Range("E112").Value = 0.1930981018275324
$ws.Range("E113").Value = 0.19782569372870323
$ws.Range("E114").Value = 0.16699472878703805
$ws.Range("E115").Value = 0.19614947844032105
$ws.Range("E116").Value = 0.1708220078874234
$ws.Range("E117").Value = 0.211523174241075
$ws.Range("E118").Value = 0.19328994063107527
$ws.Range("E119").Value = 0.21731537653220406
$ws.Range("E120").Value = 0.21077006448261207
$ws.Range("E121").Value = 0.1982862967966156
$ws.Range("E122").Value = 0.1533817713118708
$ws.Range("E123").Value = 0.13606784444360143
